$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15-22 shift down to 16-23.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new review data.
$ws.Cells.Item(15, 1).Value = 5
$ws.Cells.Item(15, 2).Value = "Excelente produto com pós venda atencioso e competente."
$ws.Cells.Item(15, 3).Value = 45954.47363185186
$ws.Cells.Item(15, 4).Value = "ZDhlNGM1NDAtZWMwMy00OGRlLWE0ZDktM2JlODM0YzJhMWYwOjU3MDE2"
